$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'knee pads 661'
    'knee pads ocr'
    'knee pad inserts for tactical pants'
    'knee pads jbm'
    'knee protectors for toddlers'
    'protec knee pads'
    'nike thermal compression pants for men'
    'under armour compression pants youth'
    'under armour compression tights for men'
    'capri pants adidas'
    'capri pants exercise'
    'capri pants nike'
    'imucci knee pads'
    'bb knee pads'
    'knee pad and helmet'
    'knee pad basketball kids'
    'knee pad buttons'
    'knee pad climbing'
    'knee pad cover'
    'knee pad cycling'
    'knee pad inserts for work pants'
    'knee pad pants women'
    'knee pad pink'
    'knee pad scooter'
    'knee pad skating'
    'knee pad wheels'
    'neoprene knee pads'
    'skating knee pads'
    'knee pads capezio'
    'knee pads canoe'
    'knee pads cycle'
    'knee pads firefighter'
    'knee pads gymnastics'
    'knee pads multicam'
    'knee pads protec'
    'knee pads shin'
    'knee pads silver'
    'skate knee pads'
    'ski knee pad'
    'kids knee protector'
    'mens warming compression pants'
    'mens workout tights'
    'womens knee pads basketball'
    'nike pro compression tights men'
    'mens basketball pants'
    'basketball knee pads for kids boys'
    'nike basketball tights'
    'military pants with knee pads'
    'asics knee pads'
    'knee pads xlarge'
    'motorcycle knee pads men'
    'gray baseball pants mens'
    'baseball pants mens knickers'
    'kids basketball knee pads youth'
    'youth knee pads basketball for kids'
    'workout leggings for men'
    'legging for men nike'
    'adidas capris men'
    'youth knee and elbow pads'
    'knee pads for dancers'
    'knee pad bathtub'
    'under armour compression tights men'
    'knee pad for dancers'
    'nike youth basketball tights'
    'compression pants women'
    'mens workout tights pants'
    'white nike compression pants men'
    'defender mens compression pants'
    'mens nike basketball pants'
    'nike basketball pants men'
    'mens leggings compression nike'
    'white compression pants men'
    'od green pants with knee pads'
    'send knee pad'
    'children knee pads'
    'smith knee pads'
    'youth xl football pants'
    'labor knee pads'
    'rubber knee pads'
    'compression pants nike'
    'athletic capris for women'
    'skins tights men'
    'supportive knee pads'
    'dye knee pads'
    'dancers knee pads'
    'apex knee pads'
    'elbow and knee pads'
    'fuse knee pads'
    'ama knee pads'
    'kp knee pads'
    'adidas tights men'
    'caterpillar knee pads'
    'husky knee pads'
    'mens leggings white'
    'nike youth compression pants'
    'muscle leggings men'
    'mens compression pants under armour'
    'airsoft pants with knee pads'
    'football pants youth with pads'
    'lotus leggings men'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
